# Auto-generated Excel COM-interop script to apply the scheduled-runner update
# to the Ridill_Profits workbook. Updates columns H-N (currentAveragePrice,
# currentAveragePriceNQ, currentAveragePriceHQ, LevePriceNQ, LevePriceHQ,
# LeveProfitNQ, LeveProfitHQ) for the affected Leve rows on each class sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 21741056
$ws.Range("I28").Value = 26318024
$ws.Range("J28").Value = 451.25
$ws.Range("K28").Value = 26318024
$ws.Range("L28").Value = 451.25
$ws.Range("M28").Value = -26317539
$ws.Range("N28").Value = -1421.25
$ws.Range("H62").Value = 92117840
$ws.Range("I62").Value = 55575444
$ws.Range("J62").Value = 125006000
$ws.Range("K62").Value = 55575444
$ws.Range("L62").Value = 125006000
$ws.Range("M62").Value = -55574820
$ws.Range("N62").Value = -125007248
$ws.Range("H65").Value = 92117840
$ws.Range("I65").Value = 55575444
$ws.Range("J65").Value = 125006000
$ws.Range("K65").Value = 277877220
$ws.Range("L65").Value = 625030000
$ws.Range("M65").Value = -277874100
$ws.Range("N65").Value = -625036240
$ws.Range("H103").Value = 35719588
$ws.Range("I103").Value = 111111640
$ws.Range("J103").Value = 7562.1055
$ws.Range("K103").Value = 333334920
$ws.Range("L103").Value = 22686.3165
$ws.Range("M103").Value = -333334334
$ws.Range("N103").Value = -23858.3165
$ws.Range("H131").Value = 4888.136
$ws.Range("I131").Value = 1031.3572
$ws.Range("J131").Value = 11637.5
$ws.Range("K131").Value = 3094.0716
$ws.Range("L131").Value = 34912.5
$ws.Range("M131").Value = 1945.9284
$ws.Range("N131").Value = -44992.5
$ws.Range("H132").Value = 2471005.8
$ws.Range("I132").Value = 1491.7675
$ws.Range("J132").Value = 55565556
$ws.Range("K132").Value = 4475.3025
$ws.Range("L132").Value = 166696668
$ws.Range("M132").Value = -1945.3025
$ws.Range("N132").Value = -166701728
$ws.Range("H137").Value = 18397776
$ws.Range("I137").Value = 1129.04
$ws.Range("J137").Value = 69499576
$ws.Range("K137").Value = 3387.12
$ws.Range("L137").Value = 208498728
$ws.Range("M137").Value = -837.1199999999999
$ws.Range("N137").Value = -208503828
$ws.Range("H141").Value = 1401.6177
$ws.Range("I141").Value = 1088.5
$ws.Range("J141").Value = 3750
$ws.Range("K141").Value = 3265.5
$ws.Range("L141").Value = 11250
$ws.Range("M141").Value = 1914.5
$ws.Range("N141").Value = -21610

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3350439.8
$ws.Range("I61").Value = 1667454.6
$ws.Range("J61").Value = 11765366
$ws.Range("K61").Value = 1667454.6
$ws.Range("L61").Value = 11765366
$ws.Range("M61").Value = -1667242.6
$ws.Range("N61").Value = -11765790
$ws.Range("H74").Value = 36509464
$ws.Range("I74").Value = 34483556
$ws.Range("J74").Value = 41028796
$ws.Range("K74").Value = 34483556
$ws.Range("L74").Value = 41028796
$ws.Range("M74").Value = -34482682
$ws.Range("N74").Value = -41030544
$ws.Range("H77").Value = 36509464
$ws.Range("I77").Value = 34483556
$ws.Range("J77").Value = 41028796
$ws.Range("K77").Value = 172417780
$ws.Range("L77").Value = 205143980
$ws.Range("M77").Value = -172413412
$ws.Range("N77").Value = -205152716
$ws.Range("H136").Value = 3350439.8
$ws.Range("I136").Value = 1667454.6
$ws.Range("J136").Value = 11765366
$ws.Range("K136").Value = 5002363.800000001
$ws.Range("L136").Value = 35296098
$ws.Range("M136").Value = -4999813.800000001
$ws.Range("N136").Value = -35301198

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 14479586
$ws.Range("I134").Value = 14706709
$ws.Range("J134").Value = 11905533
$ws.Range("K134").Value = 44120127
$ws.Range("L134").Value = 35716599
$ws.Range("M134").Value = -44117592
$ws.Range("N134").Value = -35721669

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1978.0625
$ws.Range("I16").Value = 1868.091
$ws.Range("J16").Value = 2220
$ws.Range("K16").Value = 1868.091
$ws.Range("L16").Value = 2220
$ws.Range("M16").Value = -1581.091
$ws.Range("N16").Value = -2794
$ws.Range("H31").Value = 1739449.8
$ws.Range("I31").Value = 1040.6552
$ws.Range("J31").Value = 8941431
$ws.Range("K31").Value = 1040.6552
$ws.Range("L31").Value = 8941431
$ws.Range("M31").Value = -745.6551999999999
$ws.Range("N31").Value = -8942021
$ws.Range("H34").Value = 1739449.8
$ws.Range("I34").Value = 1040.6552
$ws.Range("J34").Value = 8941431
$ws.Range("K34").Value = 1040.6552
$ws.Range("L34").Value = 8941431
$ws.Range("M34").Value = -838.6551999999999
$ws.Range("N34").Value = -8941835
$ws.Range("H50").Value = 13303.375
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 13303.375
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 13303.375
$ws.Range("N50").Value = -14553.375
$ws.Range("H58").Value = 786065.75
$ws.Range("I58").Value = 2620.8774
$ws.Range("J58").Value = 5051488
$ws.Range("K58").Value = 2620.8774
$ws.Range("L58").Value = 5051488
$ws.Range("M58").Value = -2417.8774
$ws.Range("N58").Value = -5051894
$ws.Range("H59").Value = 16159.5
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 16159.5
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 16159.5
$ws.Range("N59").Value = -18449.5
$ws.Range("H60").Value = 6888.25
$ws.Range("I60").Value = 3000
$ws.Range("J60").Value = 8184.3335
$ws.Range("K60").Value = 3000
$ws.Range("L60").Value = 8184.3335
$ws.Range("M60").Value = -2489
$ws.Range("N60").Value = -9206.333500000001
$ws.Range("H113").Value = 1978.0625
$ws.Range("I113").Value = 1868.091
$ws.Range("J113").Value = 2220
$ws.Range("K113").Value = 1868.091
$ws.Range("L113").Value = 2220
$ws.Range("M113").Value = 301.9090000000001
$ws.Range("N113").Value = -6560
$ws.Range("H134").Value = 976761.2
$ws.Range("I134").Value = 1159.2051
$ws.Range("J134").Value = 20001000
$ws.Range("K134").Value = 3477.615299999999
$ws.Range("L134").Value = 60003000
$ws.Range("M134").Value = -942.6152999999995
$ws.Range("N134").Value = -60008070
$ws.Range("H136").Value = 786065.75
$ws.Range("I136").Value = 2620.8774
$ws.Range("J136").Value = 5051488
$ws.Range("K136").Value = 7862.6322
$ws.Range("L136").Value = 15154464
$ws.Range("M136").Value = -5312.6322
$ws.Range("N136").Value = -15159564

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1954.7142
$ws.Range("I39").Value = 350
$ws.Range("J39").Value = 2392.3635
$ws.Range("K39").Value = 1050
$ws.Range("L39").Value = 7177.0905
$ws.Range("M39").Value = -756
$ws.Range("N39").Value = -7765.0905
$ws.Range("H57").Value = 2340
$ws.Range("I57").Value = 333.33334
$ws.Range("J57").Value = 3845
$ws.Range("K57").Value = 1000.00002
$ws.Range("L57").Value = 11535
$ws.Range("M57").Value = -441.0000200000001
$ws.Range("N57").Value = -12653

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 20667.166
$ws.Range("I80").Value = 9750
$ws.Range("J80").Value = 26125.75
$ws.Range("K80").Value = 9750
$ws.Range("L80").Value = 26125.75
$ws.Range("M80").Value = -8752
$ws.Range("N80").Value = -28121.75
$ws.Range("H83").Value = 20667.166
$ws.Range("I83").Value = 9750
$ws.Range("J83").Value = 26125.75
$ws.Range("K83").Value = 48750
$ws.Range("L83").Value = 130628.75
$ws.Range("M83").Value = -43758
$ws.Range("N83").Value = -140612.75
$ws.Range("H97").Value = 10000767
$ws.Range("I97").Value = 692
$ws.Range("J97").Value = 31250926
$ws.Range("K97").Value = 692
$ws.Range("L97").Value = 31250926
$ws.Range("M97").Value = -196
$ws.Range("N97").Value = -31251918
$ws.Range("H126").Value = 11121.5
$ws.Range("I126").Value = 17516.834
$ws.Range("J126").Value = 1528.5
$ws.Range("K126").Value = 52550.50199999999
$ws.Range("L126").Value = 4585.5
$ws.Range("M126").Value = -50080.50199999999
$ws.Range("N126").Value = -9525.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 971.3570999999999
$ws.Range("I61").Value = 819.4
$ws.Range("J61").Value = 1351.25
$ws.Range("K61").Value = 819.4
$ws.Range("L61").Value = 1351.25
$ws.Range("M61").Value = -617.4
$ws.Range("N61").Value = -1755.25
$ws.Range("H113").Value = 971.3570999999999
$ws.Range("I113").Value = 819.4
$ws.Range("J113").Value = 1351.25
$ws.Range("K113").Value = 819.4
$ws.Range("L113").Value = 1351.25
$ws.Range("M113").Value = 1350.6
$ws.Range("N113").Value = -5691.25
$ws.Range("H132").Value = 2307377
$ws.Range("I132").Value = 2508811.8
$ws.Range("J132").Value = 11020
$ws.Range("K132").Value = 7526435.399999999
$ws.Range("L132").Value = 33060
$ws.Range("M132").Value = -7523905.399999999
$ws.Range("N132").Value = -38120

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 386.3158
$ws.Range("I113").Value = 338.46155
$ws.Range("J113").Value = 490
$ws.Range("K113").Value = 1015.38465
$ws.Range("L113").Value = 1470
$ws.Range("M113").Value = 1154.61535
$ws.Range("N113").Value = -5810
